$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late"), shifting
# "Late" / "heading" / "Outstanding" one column to the right.
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet / tab, with the new selection.
$ws.Activate() | Out-Null
$ws.Range("R9").Select() | Out-Null
